$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$data = @(
    @(16, "79535808", "JESUS FERNANDO CASTAÑO DIAZ", "1608", 59658, 2420000),
    @(17, "1143329912", "ELIANA MARGARITA CANO ZARATE", "1608", 41600, 1040000),
    @(18, "72257127", "ALEXANDER JIMENO ALBA", "1610", 4891, 733621),
    @(19, "79535808", "JESUS FERNANDO CASTAÑO DIAZ", "1610", 59658, 2420000),
    @(20, "1143329912", "ELIANA MARGARITA CANO ZARATE", "1610", 41600, 1040000),
    @(21, "79535808", "JESUS FERNANDO CASTAÑO DIAZ", "1611", 59658, 2420000),
    @(22, "1143329912", "ELIANA MARGARITA CANO ZARATE", "1611", 41600, 1040000),
    @(23, "79535808", "JESUS FERNANDO CASTAÑO DIAZ", "1702", 59658, 2420000),
    @(24, "1143329912", "ELIANA MARGARITA CANO ZARATE", "1702", 41600, 1040000),
    @(25, "79535808", "JESUS FERNANDO CASTAÑO DIAZ", "1703", 59658, 2420000),
    @(26, "1143329912", "ELIANA MARGARITA CANO ZARATE", "1703", 41600, 1040000),
    @(27, "79535808", "JESUS FERNANDO CASTAÑO DIAZ", "1705", 59658, 2420000),
    @(28, "1143329912", "ELIANA MARGARITA CANO ZARATE", "1705", 41600, 1040000),
    @(29, "79535808", "JESUS FERNANDO CASTAÑO DIAZ", "1706", 59658, 2420000),
    @(30, "1143329912", "ELIANA MARGARITA CANO ZARATE", "1706", 41600, 1040000),
    @(31, "79535808", "JESUS FERNANDO CASTAÑO DIAZ", "1707", 59658, 2420000),
    @(32, "1143329912", "ELIANA MARGARITA CANO ZARATE", "1707", 41600, 1040000),
    @(33, "79535808", "JESUS FERNANDO CASTAÑO DIAZ", "1708", 59658, 2420000),
    @(34, "1143329912", "ELIANA MARGARITA CANO ZARATE", "1708", 41600, 1040000),
    @(35, "79535808", "JESUS FERNANDO CASTAÑO DIAZ", "1709", 59658, 2420000),
    @(36, "1143329912", "ELIANA MARGARITA CANO ZARATE", "1709", 41600, 1040000),
    @(37, "79535808", "JESUS FERNANDO CASTAÑO DIAZ", "1710", 59658, 2420000),
    @(38, "1143329912", "ELIANA MARGARITA CANO ZARATE", "1710", 41600, 1040000),
    @(39, "79535808", "JESUS FERNANDO CASTAÑO DIAZ", "1711", 59658, 2420000),
    @(40, "1143329912", "ELIANA MARGARITA CANO ZARATE", "1711", 41600, 1040000),
    @(41, "79535808", "JESUS FERNANDO CASTAÑO DIAZ", "1712", 59658, 2420000),
    @(42, "1143329912", "ELIANA MARGARITA CANO ZARATE", "1712", 41600, 1040000),
    @(43, "79535808", "JESUS FERNANDO CASTAÑO DIAZ", "1801", 59658, 2420000),
    @(44, "1143329912", "ELIANA MARGARITA CANO ZARATE", "1801", 41600, 1040000),
    @(45, "79535808", "JESUS FERNANDO CASTAÑO DIAZ", "1802", 59658, 2420000),
    @(46, "1143329912", "ELIANA MARGARITA CANO ZARATE", "1802", 41600, 1040000),
    @(47, "79535808", "JESUS FERNANDO CASTAÑO DIAZ", "1803", 59658, 2420000),
    @(48, "1143329912", "ELIANA MARGARITA CANO ZARATE", "1803", 41600, 1040000),
    @(49, "79535808", "JESUS FERNANDO CASTAÑO DIAZ", "1804", 59658, 2420000),
    @(50, "1143329912", "ELIANA MARGARITA CANO ZARATE", "1804", 41600, 1040000),
    @(51, "79535808", "JESUS FERNANDO CASTAÑO DIAZ", "1805", 59658, 2420000),
    @(52, "1143329912", "ELIANA MARGARITA CANO ZARATE", "1805", 41600, 1040000),
    @(53, "79535808", "JESUS FERNANDO CASTAÑO DIAZ", "1806", 59658, 2420000),
    @(54, "1143329912", "ELIANA MARGARITA CANO ZARATE", "1806", 41600, 1040000),
    @(55, "79535808", "JESUS FERNANDO CASTAÑO DIAZ", "1807", 59658, 2420000),
    @(56, "1143329912", "ELIANA MARGARITA CANO ZARATE", "1807", 41600, 1040000),
    @(57, "79535808", "JESUS FERNANDO CASTAÑO DIAZ", "1808", 59658, 2420000),
    @(58, "1143329912", "ELIANA MARGARITA CANO ZARATE", "1808", 41600, 1040000),
    @(59, "79535808", "JESUS FERNANDO CASTAÑO DIAZ", "1809", 59658, 2420000),
    @(60, "1143329912", "ELIANA MARGARITA CANO ZARATE", "1809", 41600, 1040000),
    @(61, "79535808", "JESUS FERNANDO CASTAÑO DIAZ", "1810", 59658, 2420000),
    @(62, "1143329912", "ELIANA MARGARITA CANO ZARATE", "1810", 41600, 1040000),
    @(63, "79535808", "JESUS FERNANDO CASTAÑO DIAZ", "1811", 59658, 2420000),
    @(64, "1143329912", "ELIANA MARGARITA CANO ZARATE", "1811", 41600, 1040000),
    @(65, "79535808", "JESUS FERNANDO CASTAÑO DIAZ", "1812", 59658, 2420000),
    @(66, "1143329912", "ELIANA MARGARITA CANO ZARATE", "1812", 41600, 1040000),
    @(67, "79535808", "JESUS FERNANDO CASTAÑO DIAZ", "1901", 59658, 2420000),
    @(68, "1143329912", "ELIANA MARGARITA CANO ZARATE", "1901", 41600, 1040000),
    @(69, "79535808", "JESUS FERNANDO CASTAÑO DIAZ", "1902", 59658, 2420000),
    @(70, "1143329912", "ELIANA MARGARITA CANO ZARATE", "1902", 41600, 1040000),
    @(71, "79535808", "JESUS FERNANDO CASTAÑO DIAZ", "1903", 59658, 2420000),
    @(72, "1143329912", "ELIANA MARGARITA CANO ZARATE", "1903", 41600, 1040000),
    @(73, "79535808", "JESUS FERNANDO CASTAÑO DIAZ", "1904", 59658, 2420000),
    @(74, "1143329912", "ELIANA MARGARITA CANO ZARATE", "1904", 41600, 1040000),
    @(75, "79535808", "JESUS FERNANDO CASTAÑO DIAZ", "1905", 59658, 2420000),
    @(76, "1143329912", "ELIANA MARGARITA CANO ZARATE", "1905", 41600, 1040000),
    @(77, "79535808", "JESUS FERNANDO CASTAÑO DIAZ", "1911", 59658, 2420000),
    @(78, "1143329912", "ELIANA MARGARITA CANO ZARATE", "1911", 41600, 1040000),
    @(79, "79535808", "JESUS FERNANDO CASTAÑO DIAZ", "1912", 59658, 2420000),
    @(80, "1143329912", "ELIANA MARGARITA CANO ZARATE", "1912", 41600, 1040000),
    @(81, "79535808", "JESUS FERNANDO CASTAÑO DIAZ", "2003", 59658, 2420000),
    @(82, "1143329912", "ELIANA MARGARITA CANO ZARATE", "2003", 41600, 1040000)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 3).Value = $row[1]
    $ws.Cells.Item($r, 4).Value = $row[2]
    $ws.Cells.Item($r, 5).Value = $row[3]
    $ws.Cells.Item($r, 6).Value = $row[4]
    $ws.Cells.Item($r, 7).Value = $row[5]
}
